$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activités")

# ---------------------------------------------------------------------------
# 1. Remove the old "total" formula in D21 and extend the table formatting
#    (copied from the last real data row, 18) down through row 31 so the new
#    entries + the now-blank trailing rows match the existing look.
# ---------------------------------------------------------------------------
$ws.Range("D21").ClearContents() | Out-Null

$ws.Range("A18:H18").Copy() | Out-Null
$ws.Range("A19:H31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. New work-log entries for "19 mai 2020".
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "19 mai 2020"
$ws.Range("B19").Value = 0.33333333333333331
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Formula = "=IF(ISBLANK(C19), NOW(),C19)-IF(ISBLANK(B19),NOW(),B19)"
$ws.Range("E19").Value = "Communication entre le client et le serveur"
$ws.Range("F19").Value = "Réalisation"
$ws.Range("G19").Value = "J'ai regardé quelques exemples de communication TCP/IP en utilisant les sockets puis j'ai adapté mon code.`nJ'ai appris à quoi servent les méthodes ""delegate"""
$ws.Hyperlinks.Add($ws.Range("H19"), "https://www.youtube.com/watch?v=Bq1JhTHlxek", "", "", "https://www.youtube.com/watch?v=Bq1JhTHlxek`nhttps://stackoverflow.com/questions/661561/how-do-i-update-the-gui-from-another-thread") | Out-Null

$ws.Range("A20").Value = "19 mai 2020"
$ws.Range("B20").Value = 0.5625
$ws.Range("C20").Value = 0.60416666666666663
$ws.Range("D20").Formula = "=IF(ISBLANK(C20), NOW(),C20)-IF(ISBLANK(B20),NOW(),B20)"
$ws.Range("E20").Value = "Communication entre le client et le serveur"
$ws.Range("F20").Value = "Réalisation"
$ws.Range("G20").Value = "Le client et le serveur peuvent maintenant communiquer de manière asynchrone."

$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 30

# ---------------------------------------------------------------------------
# 3. Extend the dropdown (list) data validations to cover the new rows.
# ---------------------------------------------------------------------------
$ws.Range("E19:E31").Validation.Add(3, 1, 1, "=Données!`$A`$2:`$A`$9") | Out-Null
$ws.Range("F19,F21:F31").Validation.Add(3, 1, 1, "=Données!`$B`$2:`$B`$6") | Out-Null

# ---------------------------------------------------------------------------
# 4. A couple of stray formatted cells below the table (left over from the
#    user's editing) - no content, just number formats.
# ---------------------------------------------------------------------------
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# 5. Selection / cursor position, matching where the user ended up.
# ---------------------------------------------------------------------------
$ws.Range("E34").Select() | Out-Null
